# ---------------------------------------------------------------------------
# Edit: B1- and B2- PowerPoint.pptx
#
# 1) Slide 5's table (the "types of financial documents" table) switches its
#    table style from the deck's custom "Table_0" style
#    ({993660BE-ADAE-43B1-85F3-BC4A3460A41C}, defined in tableStyles.xml) to
#    the built-in style {8424FB87-4E55-44C2-9CD9-30B52B7823BC}.
#
# 2) The presentation's theme switches from the "Integral" / "Red Violet"
#    colour scheme to the standard "Office" colour scheme (the deck's
#    design goes from the pink/purple Integral look to the default blue
#    Office look).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5, shape 2 (the 3-column table) --------------
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$sh.Table.ApplyStyle("{8424FB87-4E55-44C2-9CD9-30B52B7823BC}")

# --- 2) Swap the active theme's colour scheme: Integral -> Office ---------
$theme = $p.SlideMaster.Theme
$cs    = $theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
